# Add a new worksheet "PSA_LOLO" after the existing "OverallRebateEfficiency"
# sheet, populate it with a small header/data table, and make it the
# active (selected) sheet in the workbook.

$wb = $excel.ActiveWorkbook

# Worksheets.Add() inserts the new sheet before the active sheet by default,
# so move it to sit after the existing sheet to get the desired tab order.
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "PSA_LOLO"
$newSheet.Move($null, $wb.Worksheets.Item("OverallRebateEfficiency"))

# Re-resolve the sheet by name after the move - the collection order changed
# underneath the original object reference.
$newSheet = $wb.Worksheets.Item("PSA_LOLO")

# Header row, using the same (9pt) font size as the headers on the other sheet.
$newSheet.Range("A1:B1").Font.Size = 9
$newSheet.Range("A1").Value = "psa_lolo_20"
$newSheet.Range("B1").Value = "psa_loll_40"

# Data row.
$newSheet.Range("A2").Value = 3900
$newSheet.Range("B2").Value = 6200

# Match the other sheet's selected cell and make this the active tab.
[void]$newSheet.Range("A3").Select()
$newSheet.Activate()
